$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Plain-decimal-looking strings (e.g. "0.9993") get silently coerced to
    # a numeric value by the COM Value setter, which would lose the exact
    # textual representation used in this sheet (and introduce float
    # rounding). Force those to stay text, then strip the Text number
    # format back off so the cell's style index is untouched.
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $rng.NumberFormat = "@"
        $rng.Value = $text
        $rng.ClearFormats()
    } else {
        $rng.Value = $text
    }
}

Set-TextValue $ws "D2" "29.852.92"
Set-TextValue $ws "E2" "  -0.27%  "
Set-TextValue $ws "D3" "1.869.36"
Set-TextValue $ws "E3" "  -1.34%  "
Set-TextValue $ws "D4" "0.9993"
Set-TextValue $ws "E4" "  -0.16%  "
Set-TextValue $ws "D5" "0.7379"
Set-TextValue $ws "E5" "  -4.80%  "
Set-TextValue $ws "D6" "241.95"
Set-TextValue $ws "E6" "  -1.15%  "
Set-TextValue $ws "D7" "0.9995"
Set-TextValue $ws "E7" "  -0.15%  "
Set-TextValue $ws "D8" "0.3150"
Set-TextValue $ws "E8" "  +0.39%  "
Set-TextValue $ws "D9" "24.71"
Set-TextValue $ws "E9" "  -4.31%  "
Set-TextValue $ws "D10" "0.07122"
Set-TextValue $ws "E10" "  -1.78%  "
Set-TextValue $ws "D11" "0.08386"
Set-TextValue $ws "E11" "  -6.11%  "
Set-TextValue $ws "D12" "0.7528"
Set-TextValue $ws "E12" "  -2.62%  "
Set-TextValue $ws "D13" "5.465"
Set-TextValue $ws "E13" "  +0.07%  "
Set-TextValue $ws "D14" "1.890.85"
Set-TextValue $ws "E14" "  -0.72%  "
Set-TextValue $ws "D15" "92.48"
Set-TextValue $ws "E15" "  -2.41%  "
Set-TextValue $ws "D16" "29.860.53"
Set-TextValue $ws "E16" "  -0.38%  "
Set-TextValue $ws "D17" "6.041"
Set-TextValue $ws "E17" "  -2.72%  "
Set-TextValue $ws "D18" "13.58"
Set-TextValue $ws "E18" "  -2.95%  "
Set-TextValue $ws "D19" "243.06"
Set-TextValue $ws "E19" "  -1.50%  "
Set-TextValue $ws "D20" "0.000007834"
Set-TextValue $ws "E20" "  -0.80%  "
Set-TextValue $ws "D21" "0.9987"
Set-TextValue $ws "E21" "  -0.22%  "
Set-TextValue $ws "D22" "2.115.15"
Set-TextValue $ws "E22" "  -3.12%  "
Set-TextValue $ws "D23" "7.911"
Set-TextValue $ws "E23" "  -3.33%  "
Set-TextValue $ws "D24" "0.9995"
Set-TextValue $ws "E24" "  -0.14%  "
Set-TextValue $ws "D25" "0.1564"
Set-TextValue $ws "E25" "  -1.79%  "
Set-TextValue $ws "D26" "9.315"
Set-TextValue $ws "E26" "  -2.67%  "
Set-TextValue $ws "D27" "164.30"
Set-TextValue $ws "E27" "  +0.90%  "
Set-TextValue $ws "D28" "18.58"
Set-TextValue $ws "E28" "  -1.46%  "
Set-TextValue $ws "D29" "2.017"
Set-TextValue $ws "E29" "  -1.48%  "
Set-TextValue $ws "D30" "1.476"
Set-TextValue $ws "E30" "  +3.47%  "
Set-TextValue $ws "D31" "4.653"
Set-TextValue $ws "E31" "  +2.76%  "
Set-TextValue $ws "D32" "1.530"
Set-TextValue $ws "E32" "  -1.32%  "
Set-TextValue $ws "D33" "4.320"
Set-TextValue $ws "E33" "  +4.90%  "
Set-TextValue $ws "D34" "0.05330"
Set-TextValue $ws "E34" "  -3.26%  "
Set-TextValue $ws "D35" "1.235"
Set-TextValue $ws "E35" "  -1.01%  "
Set-TextValue $ws "D36" "0.7538"
Set-TextValue $ws "E36" "  -0.12%  "
Set-TextValue $ws "E37" "  +0.04%  "
Set-TextValue $ws "D38" "2.700"
Set-TextValue $ws "E38" "  -0.42%  "
Set-TextValue $ws "D39" "0.01956"
Set-TextValue $ws "E39" "  -0.55%  "
Set-TextValue $ws "D40" "2.747"
Set-TextValue $ws "E40" "  -1.61%  "
Set-TextValue $ws "D41" "0.4476"
Set-TextValue $ws "E41" "  -0.77%  "
Set-TextValue $ws "D42" "1.099.98"
Set-TextValue $ws "E42" "  +0.67%  "
Set-TextValue $ws "D43" "6.073"
Set-TextValue $ws "E43" "  -0.21%  "
Set-TextValue $ws "D44" "72.28"
Set-TextValue $ws "E44" "  -2.69%  "
Set-TextValue $ws "D45" "0.8614"
Set-TextValue $ws "E45" "  +0.66%  "
Set-TextValue $ws "D46" "1.000"
Set-TextValue $ws "E46" "  -0.04%  "
Set-TextValue $ws "D47" "103.09"
Set-TextValue $ws "E47" "  +0.19%  "
Set-TextValue $ws "D48" "7.696"
Set-TextValue $ws "E48" "  +0.71%  "
Set-TextValue $ws "E49" "  -2.86%  "
Set-TextValue $ws "D50" "3.056"
Set-TextValue $ws "E50" "  +2.16%  "
Set-TextValue $ws "D51" "2.015.47"
Set-TextValue $ws "E51" "  -2.06%  "
